# Apply resume content updates via Find/Replace across the document.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $find
    }
}

# 1. Skills line: "Node.js" -> "Node, React"
Replace-Text "Python, JavaScript, HTML, CSS, Node.js, TypeScript, Java, C++, SQL" `
             "Python, JavaScript, HTML, CSS, Node, React, TypeScript, Java, C++, SQL"

# 2. Tools line: reorder items
Replace-Text "Git, macOS, Windows, Visual Studio, Eclipse, Unix, Linux" `
             "Git, Unix, Linux, Visual Studio, Eclipse, macOS, Windows"

# 3. "Built 8 applications..." bullet
Replace-Text "Built 8 applications end-to-end and worked in 3 other applications." `
             "Built 8 applications end-to-end and worked in other applications assisting in development."

# 4. "Meta Tag Generator:" -> "Metadata Generator:"
Replace-Text "Meta Tag Generator: " "Metadata Generator: "

# 5. Plugin description rewrite
Replace-Text "Developed an Eleventy plugin with JavaScript that has 170+ users. Generate document metadata for the <head> of a webpage containing: Open Graph, Twitter card, generic meta tags and a canonical link." `
             "Developed an Eleventy plugin with JavaScript that generates document metadata for the <head> of a webpage containing: Open Graph, Twitter card, generic meta tags, CSS, JS, custom tags, and a canonical link. Plugin is published on npm and currently has 180+ users."

# 6. Eleventy Photo Gallery bullet
Replace-Text "Created a responsive image gallery site template using the Eleventy static site generator. Images are dynamically generated with Node.js at build time." `
             "Created a responsive image gallery site template using the Eleventy static site generator. Responsive images are dynamically generated with Node.js at build time."
